$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove two data rows entirely (error-calc cleanup) ---
# Delete higher row index first so the lower one's index stays valid.
$ws.Rows.Item(28).Delete()   # "SC 92" row removed
$ws.Rows.Item(26).Delete()   # "RM 232" row removed

# --- Column D (sheet column E) value corrections on the remaining top block ---
$ws.Range("E2").Value = -7.2
$ws.Range("E6").ClearContents()
$ws.Range("E12").Value = -5.3
$ws.Range("E14").ClearContents()
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E23").ClearContents()
$ws.Range("E24").ClearContents()

# --- Corrections on the shifted SC block (rows 26-33 after the deletions above) ---
$ws.Range("C26").Value = 10.8          # SC 5
$ws.Range("C27").ClearContents()       # SC 101
$ws.Range("C30").Value = 11.4          # SC 120
$ws.Range("E31").Value = -8.1          # SC 132
$ws.Range("C32").ClearContents()       # SC 193
$ws.Range("E33").Value = -10.7         # SC 232
